$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct Timezone/Timezone_path test data for GMT-12 row (row 2) and shift rows 50-76 ---
# Row 2 Timezone_path ("G2") was incorrectly duplicated as "India Standard Time";
# correct value pairs with H2 "(GMT-12:00) International Date Line West" is "Dateline Standard Time".
$ws.Range("G2").Value2 = "Dateline Standard Time"
$ws.Range("H2").Value2 = "(GMT-12:00) International Date Line West"
$ws.Range("G50").Value2 = "Nepal Standard Time"
$ws.Range("H50").Value2 = "(GMT+05:45) Kathmandu"
$ws.Range("G51").Value2 = "Central Asia Standard Time"
$ws.Range("H51").Value2 = "(GMT+06:00) Astana, Dhaka, Almaty"
$ws.Range("G52").Value2 = "Sri Lanka Standard Time"
$ws.Range("H52").Value2 = "(GMT+05:30) Sri Jayawardenepura"
$ws.Range("G53").Value2 = "N. Central Asia Standard Time"
$ws.Range("H53").Value2 = "(GMT+07:00) Novosibirsk"
$ws.Range("G54").Value2 = "Myanmar Standard Time"
$ws.Range("H54").Value2 = "(GMT+06:30) Rangoon"
$ws.Range("G55").Value2 = "SE Asia Standard Time"
$ws.Range("H55").Value2 = "(GMT+07:00) Bangkok, Hanoi, Jakarta"
$ws.Range("G56").Value2 = "North Asia Standard Time"
$ws.Range("H56").Value2 = "(GMT+07:00) Krasnoyarsk"
$ws.Range("G57").Value2 = "China Standard Time"
$ws.Range("H57").Value2 = "(GMT+08:00) Beijing, Chongqing, Hong Kong, Urumqi"
$ws.Range("G58").Value2 = "Singapore Standard Time"
$ws.Range("H58").Value2 = "(GMT+08:00) Kuala Lumpur, Singapore"
$ws.Range("G59").Value2 = "Taipei Standard Time"
$ws.Range("H59").Value2 = "(GMT+08:00) Taipei"
$ws.Range("G60").Value2 = "W. Australia Standard Time"
$ws.Range("H60").Value2 = "(GMT+08:00) Perth"
$ws.Range("G61").Value2 = "North Asia East Standard Time"
$ws.Range("H61").Value2 = "(GMT+08:00) Irkutsk, Ulaan Bataar"
$ws.Range("G62").Value2 = "Korea Standard Time"
$ws.Range("H62").Value2 = "(GMT+09:00) Seoul"
$ws.Range("G63").Value2 = "Tokyo Standard Time"
$ws.Range("H63").Value2 = "(GMT+09:00) Osaka, Sapporo, Tokyo"
$ws.Range("G64").Value2 = "Yakutsk Standard Time"
$ws.Range("H64").Value2 = "(GMT+09:00) Yakutsk"
$ws.Range("G65").Value2 = "AUS Central Standard Time"
$ws.Range("H65").Value2 = "(GMT+09:30) Darwin"
$ws.Range("G66").Value2 = "Cen. Australia Standard Time"
$ws.Range("H66").Value2 = "(GMT+09:30) Adelaide"
$ws.Range("G67").Value2 = "AUS Eastern Standard Time"
$ws.Range("H67").Value2 = "(GMT+10:00) Canberra, Melbourne, Sydney"
$ws.Range("G68").Value2 = "E. Australia Standard Time"
$ws.Range("H68").Value2 = "(GMT+10:00) Brisbane"
$ws.Range("G69").Value2 = "Tasmania Standard Time"
$ws.Range("H69").Value2 = "(GMT+10:00) Hobart"
$ws.Range("G70").Value2 = "Vladivostok Standard Time"
$ws.Range("H70").Value2 = "(GMT+10:00) Vladivostok"
$ws.Range("G71").Value2 = "West Pacific Standard Time"
$ws.Range("H71").Value2 = "(GMT+10:00) Guam, Port Moresby"
$ws.Range("G72").Value2 = "Central Pacific Standard Time"
$ws.Range("H72").Value2 = "(GMT+11:00) Magadan, Solomon Is., New Caledonia"
$ws.Range("G73").Value2 = "Fiji Standard Time"
$ws.Range("H73").Value2 = "(GMT+12:00) Fiji"
$ws.Range("G74").Value2 = "New Zealand Standard Time"
$ws.Range("H74").Value2 = "(GMT+12:00) Auckland, Wellington"
$ws.Range("G75").Value2 = "Tonga Standard Time"
$ws.Range("H75").Value2 = "(GMT+13:00) Nuku'alofa"
$ws.Range("G76").Value2 = "India Standard Time"
$ws.Range("H76").Value2 = "(GMT+05:30) Chennai, Kolkata, Mumbai, New Delhi"

# --- Shrink the hidden AutoFilter range (table actually has 75 data-driving rows now) ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$A`$75"
    }
}

# --- Update the selection shown when the workbook is reopened ---
$ws.Rows("50:50").Select()
